$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "511÷3=" "469÷2="
Replace-Text "396÷4=" "865÷9="
Replace-Text "746÷8=" "813÷6="
Replace-Text "504÷7=" "282÷3="
Replace-Text "307÷5=" "980÷3="
Replace-Text "613÷6=" "887÷3="
Replace-Text "132÷5=" "804÷3="
Replace-Text "496÷7=" "160÷9="
Replace-Text "549÷3=" "669÷6="
Replace-Text "486÷8=" "290÷6="
Replace-Text "452÷9=" "716÷9="
Replace-Text "744÷4=" "694÷2="
Replace-Text "707÷2=" "958÷7="
Replace-Text "422÷5=" "587÷3="
Replace-Text "482÷8=" "535÷4="
Replace-Text "596÷9=" "298÷4="
Replace-Text "617÷9=" "126÷8="
Replace-Text "249÷4=" "519÷6="
Replace-Text "199÷4=" "262÷9="
Replace-Text "486÷2=" "209÷3="
Replace-Text "105÷7=" "606÷8="
Replace-Text "619÷3=" "155÷7="
Replace-Text "329÷4=" "298÷9="
Replace-Text "296÷3=" "440÷7="
Replace-Text "214÷9=" "931÷6="

Write-Host "Done"
